# Commit: "Truncate all measurements to integers in fake_data.xlsx"
#
# Round (truncate) Temperature, Variable_1 and Variable_2 down to whole
# numbers. These live in columns B, C and D (column A is the Date column
# and is left untouched, aside from re-affirming its YYYY-MM-DD display
# format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

# Data starts on row 2 (row 1 holds the column headers).
$dataFirstRow = $firstRow + 1

# Measurement columns: B = Temperature (F), C = Variable_1, D = Variable_2
$firstCol = 2
$lastCol = 4

for ($r = $dataFirstRow; $r -le $lastRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $orig = $cell.Value2
        if ($orig -ne $null) {
            # Truncate toward zero (e.g. 18.5 -> 18, 9.9 -> 9), matching the
            # integer truncation applied to the source data / HTML preview.
            $cell.Value2 = [int]$orig
        }
    }
}

# Re-apply the date display format to the Date column so it keeps rendering
# as YYYY-MM-DD.
$ws.Range("A$dataFirstRow`:A$lastRow").NumberFormat = "YYYY-MM-DD"

Write-Host "Truncated measurement columns B:D (rows $dataFirstRow-$lastRow) to integers"
